$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from A1 (bold, centered, bordered) to new header cells G1:H1
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

$ws.Range("G1").Value = "num_samples"
$ws.Range("H1").Value = "fractional_uncertainty"

$ws.Range("G2").Value = 956
$ws.Range("H2").Value = 0.03474455528897709
$ws.Range("G3").Value = 731
$ws.Range("H3").Value = 0.03492634042692064
$ws.Range("G4").Value = 971
$ws.Range("H4").Value = 0.03169069461024138
$ws.Range("G5").Value = 772
$ws.Range("H5").Value = 0.032809869967679
$ws.Range("G6").Value = 973
$ws.Range("H6").Value = 0.03009467527024692
$ws.Range("G7").Value = 783
$ws.Range("H7").Value = 0.03067124112369637
$ws.Range("G8").Value = 940
$ws.Range("H8").Value = 0.02813431984812085
$ws.Range("G9").Value = 792
$ws.Range("H9").Value = 0.03630763043839624
$ws.Range("G10").Value = 937
$ws.Range("H10").Value = 0.02850992436117416
$ws.Range("G11").Value = 804
$ws.Range("H11").Value = 0.03292225038833987
$ws.Range("G12").Value = 953
$ws.Range("H12").Value = 0.02999691922471324
$ws.Range("G13").Value = 808
$ws.Range("H13").Value = 0.03267920594117515
$ws.Range("G14").Value = 956
$ws.Range("H14").Value = 0.03148952966005882
$ws.Range("G15").Value = 796
$ws.Range("H15").Value = 0.03016095376771777
$ws.Range("G16").Value = 964
$ws.Range("H16").Value = 0.03184115097705828
$ws.Range("G17").Value = 790
$ws.Range("H17").Value = 0.02937649678565651

Write-Output "done"
